# Update computed pricing/profit columns (H:N) across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed
# currentAveragePrice / LevePrice / LeveProfit figures pulled by the
# scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4170975  # H74: was 4766348
$ws.Cells.Item(74, 9).Value = 9094400  # I74: was 12503538
$ws.Cells.Item(74, 11).Value = 9094400  # K74: was 12503538
$ws.Cells.Item(74, 13).Value = -9093464  # M74: was -12502602
$ws.Cells.Item(77, 8).Value = 4170975  # H77: was 4766348
$ws.Cells.Item(77, 9).Value = 9094400  # I77: was 12503538
$ws.Cells.Item(77, 11).Value = 45472000  # K77: was 62517690
$ws.Cells.Item(77, 13).Value = -45467320  # M77: was -62513010
$ws.Cells.Item(98, 8).Value = 3688.7646  # H98: was 3035.2
$ws.Cells.Item(98, 9).Value = 1835.7273  # I98: was 1429.2
$ws.Cells.Item(98, 10).Value = 7086  # J98: was 7853.2
$ws.Cells.Item(98, 11).Value = 1835.7273  # K98: was 1429.2
$ws.Cells.Item(98, 12).Value = 7086  # L98: was 7853.2
$ws.Cells.Item(98, 13).Value = -337.7273  # M98: was 68.79999999999995
$ws.Cells.Item(98, 14).Value = -10082  # N98: was -10849.2
$ws.Cells.Item(122, 8).Value = 3688.7646  # H122: was 3035.2
$ws.Cells.Item(122, 9).Value = 1835.7273  # I122: was 1429.2
$ws.Cells.Item(122, 10).Value = 7086  # J122: was 7853.2
$ws.Cells.Item(122, 11).Value = 5507.1819  # K122: was 4287.6
$ws.Cells.Item(122, 12).Value = 21258  # L122: was 23559.6
$ws.Cells.Item(122, 13).Value = -3057.1819  # M122: was -1837.6
$ws.Cells.Item(122, 14).Value = -26158  # N122: was -28459.6
$ws.Cells.Item(129, 8).Value = 874.7093  # H129: was 843.2033699999999
$ws.Cells.Item(129, 9).Value = 356.33334  # I129: was 359.58334
$ws.Cells.Item(129, 10).Value = 958.77026  # J129: was 966.68085
$ws.Cells.Item(129, 11).Value = 1069.00002  # K129: was 1078.75002
$ws.Cells.Item(129, 12).Value = 2876.31078  # L129: was 2900.04255
$ws.Cells.Item(129, 13).Value = 3930.99998  # M129: was 3921.24998
$ws.Cells.Item(129, 14).Value = -12876.31078  # N129: was -12900.04255
$ws.Cells.Item(137, 8).Value = 2034.5614  # H137: was 2294.3264
$ws.Cells.Item(137, 9).Value = 1044.0513  # I137: was 1195.258
$ws.Cells.Item(137, 10).Value = 4180.6665  # J137: was 4187.1665
$ws.Cells.Item(137, 11).Value = 3132.1539  # K137: was 3585.774
$ws.Cells.Item(137, 12).Value = 12541.9995  # L137: was 12561.4995
$ws.Cells.Item(137, 13).Value = -582.1539000000002  # M137: was -1035.774
$ws.Cells.Item(137, 14).Value = -17641.9995  # N137: was -17661.4995

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1149.75  # H2: was 656.85
$ws.Cells.Item(2, 9).Value = 1199.6666  # I2: was 664.25
$ws.Cells.Item(2, 10).Value = 1000  # J2: was 627.25
$ws.Cells.Item(2, 11).Value = 1199.6666  # K2: was 664.25
$ws.Cells.Item(2, 12).Value = 1000  # L2: was 627.25
$ws.Cells.Item(2, 13).Value = -1086.6666  # M2: was -551.25
$ws.Cells.Item(2, 14).Value = -1226  # N2: was -853.25
$ws.Cells.Item(32, 8).Value = 3527.5513  # H32: was 3857.2676
$ws.Cells.Item(32, 9).Value = 3169.324  # I32: was 3495.9219
$ws.Cells.Item(32, 11).Value = 3169.324  # K32: was 3495.9219
$ws.Cells.Item(32, 13).Value = -2882.324  # M32: was -3208.9219
$ws.Cells.Item(45, 8).Value = 1808.6666  # H45: was 1566.1428
$ws.Cells.Item(45, 9).Value = 1770.4  # I45: was 1493.8334
$ws.Cells.Item(45, 11).Value = 1770.4  # K45: was 1493.8334
$ws.Cells.Item(45, 13).Value = -1393.4  # M45: was -1116.8334
$ws.Cells.Item(61, 8).Value = 926.03125  # H61: was 941.9643
$ws.Cells.Item(61, 9).Value = 706  # I61: was 656.05
$ws.Cells.Item(61, 10).Value = 1410.1  # J61: was 1656.75
$ws.Cells.Item(61, 11).Value = 706  # K61: was 656.05
$ws.Cells.Item(61, 12).Value = 1410.1  # L61: was 1656.75
$ws.Cells.Item(61, 13).Value = -494  # M61: was -444.05
$ws.Cells.Item(61, 14).Value = -1834.1  # N61: was -2080.75
$ws.Cells.Item(109, 8).Value = 25560.953  # H109: was 25571.428
$ws.Cells.Item(109, 10).Value = 25560.953  # J109: was 25571.428
$ws.Cells.Item(109, 12).Value = 25560.953  # L109: was 25571.428
$ws.Cells.Item(109, 14).Value = -28334.953  # N109: was -28345.428
$ws.Cells.Item(116, 8).Value = 1149.75  # H116: was 656.85
$ws.Cells.Item(116, 9).Value = 1199.6666  # I116: was 664.25
$ws.Cells.Item(116, 10).Value = 1000  # J116: was 627.25
$ws.Cells.Item(116, 11).Value = 1199.6666  # K116: was 664.25
$ws.Cells.Item(116, 12).Value = 1000  # L116: was 627.25
$ws.Cells.Item(116, 13).Value = 1094.3334  # M116: was 1629.75
$ws.Cells.Item(116, 14).Value = -5588  # N116: was -5215.25
$ws.Cells.Item(122, 8).Value = 1975.3889  # H122: was 2228.2
$ws.Cells.Item(122, 9).Value = 1412.68  # I122: was 1608
$ws.Cells.Item(122, 10).Value = 3254.2727  # J122: was 3468.6
$ws.Cells.Item(122, 11).Value = 4238.04  # K122: was 4824
$ws.Cells.Item(122, 12).Value = 9762.8181  # L122: was 10405.8
$ws.Cells.Item(122, 13).Value = -1788.04  # M122: was -2374
$ws.Cells.Item(122, 14).Value = -14662.8181  # N122: was -15305.8
$ws.Cells.Item(132, 8).Value = 2217.6538  # H132: was 2355.6597
$ws.Cells.Item(132, 9).Value = 1513.3334  # I132: was 1587.125
$ws.Cells.Item(132, 10).Value = 3802.375  # J132: was 3995.2
$ws.Cells.Item(132, 11).Value = 4540.0002  # K132: was 4761.375
$ws.Cells.Item(132, 12).Value = 11407.125  # L132: was 11985.6
$ws.Cells.Item(132, 13).Value = -2010.0002  # M132: was -2231.375
$ws.Cells.Item(132, 14).Value = -16467.125  # N132: was -17045.6
$ws.Cells.Item(136, 8).Value = 926.03125  # H136: was 941.9643
$ws.Cells.Item(136, 9).Value = 706  # I136: was 656.05
$ws.Cells.Item(136, 10).Value = 1410.1  # J136: was 1656.75
$ws.Cells.Item(136, 11).Value = 2118  # K136: was 1968.15
$ws.Cells.Item(136, 12).Value = 4230.299999999999  # L136: was 4970.25
$ws.Cells.Item(136, 13).Value = 432  # M136: was 581.8500000000001
$ws.Cells.Item(136, 14).Value = -9330.299999999999  # N136: was -10070.25
$ws.Cells.Item(137, 8).Value = 39766  # H137: was 39766.855
$ws.Cells.Item(137, 10).Value = 39766  # J137: was 39766.855
$ws.Cells.Item(137, 12).Value = 39766  # L137: was 39766.855
$ws.Cells.Item(137, 14).Value = -49966  # N137: was -49966.855

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1149.75  # H3: was 656.85
$ws.Cells.Item(3, 9).Value = 1199.6666  # I3: was 664.25
$ws.Cells.Item(3, 10).Value = 1000  # J3: was 627.25
$ws.Cells.Item(3, 11).Value = 1199.6666  # K3: was 664.25
$ws.Cells.Item(3, 12).Value = 1000  # L3: was 627.25
$ws.Cells.Item(3, 13).Value = -1085.6666  # M3: was -550.25
$ws.Cells.Item(3, 14).Value = -1228  # N3: was -855.25
$ws.Cells.Item(59, 8).Value = 27875  # H59: was 24583
$ws.Cells.Item(59, 10).Value = 27875  # J59: was 24583
$ws.Cells.Item(59, 12).Value = 27875  # L59: was 24583
$ws.Cells.Item(59, 14).Value = -29569  # N59: was -26277
$ws.Cells.Item(80, 9).Value = 95.666664  # I80: was 83.72727
$ws.Cells.Item(80, 10).Value = 155  # J80: was 170.625
$ws.Cells.Item(80, 11).Value = 95.666664  # K80: was 83.72727
$ws.Cells.Item(80, 12).Value = 155  # L80: was 170.625
$ws.Cells.Item(80, 13).Value = 902.333336  # M80: was 914.27273
$ws.Cells.Item(80, 14).Value = -2151  # N80: was -2166.625
$ws.Cells.Item(83, 9).Value = 95.666664  # I83: was 83.72727
$ws.Cells.Item(83, 10).Value = 155  # J83: was 170.625
$ws.Cells.Item(83, 11).Value = 478.33332  # K83: was 418.63635
$ws.Cells.Item(83, 12).Value = 775  # L83: was 853.125
$ws.Cells.Item(83, 13).Value = 4513.66668  # M83: was 4573.36365
$ws.Cells.Item(83, 14).Value = -10759  # N83: was -10837.125
$ws.Cells.Item(94, 8).Value = 442.87878  # H94: was 433.05884
$ws.Cells.Item(94, 9).Value = 556.36365  # I94: was 536.913
$ws.Cells.Item(94, 11).Value = 556.36365  # K94: was 536.913
$ws.Cells.Item(94, 13).Value = -105.36365  # M94: was -85.91300000000001
$ws.Cells.Item(107, 8).Value = 1636.3334  # H107: was 1917.75
$ws.Cells.Item(107, 9).Value = 1558.6666  # I107: was 1946.2
$ws.Cells.Item(107, 10).Value = 1869.3334  # J107: was 1870.3334
$ws.Cells.Item(107, 11).Value = 1558.6666  # K107: was 1946.2
$ws.Cells.Item(107, 12).Value = 1869.3334  # L107: was 1870.3334
$ws.Cells.Item(107, 13).Value = 361.3334  # M107: was -26.20000000000005
$ws.Cells.Item(107, 14).Value = -5709.3334  # N107: was -5710.3334
$ws.Cells.Item(134, 8).Value = 2032.541  # H134: was 2092.4822
$ws.Cells.Item(134, 9).Value = 1156.4222  # I134: was 1172.721
$ws.Cells.Item(134, 10).Value = 4496.625  # J134: was 5134.769
$ws.Cells.Item(134, 11).Value = 3469.2666  # K134: was 3518.163
$ws.Cells.Item(134, 12).Value = 13489.875  # L134: was 15404.307
$ws.Cells.Item(134, 13).Value = -934.2665999999999  # M134: was -983.163
$ws.Cells.Item(134, 14).Value = -18559.875  # N134: was -20474.307
$ws.Cells.Item(137, 8).Value = 37236.668  # H137: was 37025.6
$ws.Cells.Item(137, 10).Value = 37236.668  # J137: was 37025.6
$ws.Cells.Item(137, 12).Value = 37236.668  # L137: was 37025.6
$ws.Cells.Item(137, 14).Value = -47436.668  # N137: was -47225.6
$ws.Cells.Item(140, 8).Value = 49030  # H140: was 58916.25
$ws.Cells.Item(140, 10).Value = 49030  # J140: was 58916.25
$ws.Cells.Item(140, 12).Value = 49030  # L140: was 58916.25
$ws.Cells.Item(140, 14).Value = -59390  # N140: was -69276.25

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(37, 8).Value = 15000  # H37: was 3000
$ws.Cells.Item(37, 9).Value = 0  # I37: was 3000
$ws.Cells.Item(37, 10).Value = 15000  # J37: was 0
$ws.Cells.Item(37, 11).Value = 0  # K37: was 3000
$ws.Cells.Item(37, 12).Value = 15000  # L37: was 0
$ws.Cells.Item(37, 13).ClearContents()  # M37: was -2893
$ws.Cells.Item(37, 14).Value = -15214  # N37: was None
$ws.Cells.Item(58, 8).Value = 1522.19  # H58: was 1045.6735
$ws.Cells.Item(58, 9).Value = 1355.4875  # I58: was 1155.8379
$ws.Cells.Item(58, 10).Value = 2189  # J58: was 706
$ws.Cells.Item(58, 11).Value = 1355.4875  # K58: was 1155.8379
$ws.Cells.Item(58, 12).Value = 2189  # L58: was 706
$ws.Cells.Item(58, 13).Value = -1152.4875  # M58: was -952.8379
$ws.Cells.Item(58, 14).Value = -2595  # N58: was -1112
$ws.Cells.Item(132, 8).Value = 4419.579  # H132: was 2534.4412
$ws.Cells.Item(132, 9).Value = 4280.7144  # I132: was 2439.84
$ws.Cells.Item(132, 10).Value = 4808.4  # J132: was 2797.2222
$ws.Cells.Item(132, 11).Value = 12842.1432  # K132: was 7319.52
$ws.Cells.Item(132, 12).Value = 14425.2  # L132: was 8391.6666
$ws.Cells.Item(132, 13).Value = -10312.1432  # M132: was -4789.52
$ws.Cells.Item(132, 14).Value = -19485.2  # N132: was -13451.6666
$ws.Cells.Item(134, 8).Value = 1902.9722  # H134: was 3263.182
$ws.Cells.Item(134, 9).Value = 968.1053000000001  # I134: was 4669.731
$ws.Cells.Item(134, 10).Value = 2947.8235  # J134: was 2002.138
$ws.Cells.Item(134, 11).Value = 2904.3159  # K134: was 14009.193
$ws.Cells.Item(134, 12).Value = 8843.470499999999  # L134: was 6006.414
$ws.Cells.Item(134, 13).Value = -369.3159000000001  # M134: was -11474.193
$ws.Cells.Item(134, 14).Value = -13913.4705  # N134: was -11076.414
$ws.Cells.Item(136, 8).Value = 1522.19  # H136: was 1045.6735
$ws.Cells.Item(136, 9).Value = 1355.4875  # I136: was 1155.8379
$ws.Cells.Item(136, 10).Value = 2189  # J136: was 706
$ws.Cells.Item(136, 11).Value = 4066.4625  # K136: was 3467.5137
$ws.Cells.Item(136, 12).Value = 6567  # L136: was 2118
$ws.Cells.Item(136, 13).Value = -1516.4625  # M136: was -917.5137
$ws.Cells.Item(136, 14).Value = -11667  # N136: was -7218

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1314.2927  # H5: was 1320.4884
$ws.Cells.Item(5, 9).Value = 324.42105  # I5: was 325.21054
$ws.Cells.Item(5, 10).Value = 2169.182  # J5: was 2108.4167
$ws.Cells.Item(5, 11).Value = 973.26315  # K5: was 975.6316199999999
$ws.Cells.Item(5, 12).Value = 6507.545999999999  # L5: was 6325.250100000001
$ws.Cells.Item(5, 13).Value = -861.26315  # M5: was -863.6316199999999
$ws.Cells.Item(5, 14).Value = -6731.545999999999  # N5: was -6549.250100000001
$ws.Cells.Item(131, 8).Value = 884.902  # H131: was 853.36707
$ws.Cells.Item(131, 9).Value = 695  # I131: was 649.375
$ws.Cells.Item(131, 10).Value = 915.11365  # J131: was 876.3521
$ws.Cells.Item(131, 11).Value = 2085  # K131: was 1948.125
$ws.Cells.Item(131, 12).Value = 2745.34095  # L131: was 2629.0563
$ws.Cells.Item(131, 13).Value = 2955  # M131: was 3091.875
$ws.Cells.Item(131, 14).Value = -12825.34095  # N131: was -12709.0563
$ws.Cells.Item(135, 8).Value = 1314.2927  # H135: was 1320.4884
$ws.Cells.Item(135, 9).Value = 324.42105  # I135: was 325.21054
$ws.Cells.Item(135, 10).Value = 2169.182  # J135: was 2108.4167
$ws.Cells.Item(135, 11).Value = 2919.78945  # K135: was 2926.89486
$ws.Cells.Item(135, 12).Value = 19522.638  # L135: was 18975.7503
$ws.Cells.Item(135, 13).Value = -384.7894499999998  # M135: was -391.8948599999999
$ws.Cells.Item(135, 14).Value = -24592.638  # N135: was -24045.7503

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2897.32  # H122: was 4287.2144
$ws.Cells.Item(122, 9).Value = 1401.1333  # I122: was 2101.25
$ws.Cells.Item(122, 10).Value = 5141.6  # J122: was 5161.6
$ws.Cells.Item(122, 11).Value = 4203.3999  # K122: was 6303.75
$ws.Cells.Item(122, 12).Value = 15424.8  # L122: was 15484.8
$ws.Cells.Item(122, 13).Value = -1753.3999  # M122: was -3853.75
$ws.Cells.Item(122, 14).Value = -20324.8  # N122: was -20384.8
$ws.Cells.Item(123, 8).Value = 10962.0625  # H123: was 10861.053
$ws.Cells.Item(123, 10).Value = 10962.0625  # J123: was 10861.053
$ws.Cells.Item(123, 12).Value = 10962.0625  # L123: was 10861.053
$ws.Cells.Item(123, 14).Value = -15862.0625  # N123: was -15761.053
$ws.Cells.Item(133, 8).Value = 38114.305  # H133: was 37898.137
$ws.Cells.Item(133, 10).Value = 38114.305  # J133: was 37898.137
$ws.Cells.Item(133, 12).Value = 38114.305  # L133: was 37898.137
$ws.Cells.Item(133, 14).Value = -48234.305  # N133: was -48018.137
$ws.Cells.Item(137, 8).Value = 42716.668  # H137: was 42667.6
$ws.Cells.Item(137, 10).Value = 42716.668  # J137: was 42667.6
$ws.Cells.Item(137, 12).Value = 42716.668  # L137: was 42667.6
$ws.Cells.Item(137, 14).Value = -52916.668  # N137: was -52867.6

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 37043676  # H40: was 32263860
$ws.Cells.Item(40, 9).Value = 76926720  # I40: was 58826340
$ws.Cells.Item(40, 11).Value = 76926720  # K40: was 58826340
$ws.Cells.Item(40, 13).Value = -76926584  # M40: was -58826204
$ws.Cells.Item(122, 8).Value = 5756.6113  # H122: was 6693.769
$ws.Cells.Item(122, 9).Value = 2989.25  # I122: was 3202.8
$ws.Cells.Item(122, 10).Value = 7970.5  # J122: was 8875.625
$ws.Cells.Item(122, 11).Value = 8967.75  # K122: was 9608.400000000001
$ws.Cells.Item(122, 12).Value = 23911.5  # L122: was 26626.875
$ws.Cells.Item(122, 13).Value = -6517.75  # M122: was -7158.400000000001
$ws.Cells.Item(122, 14).Value = -28811.5  # N122: was -31526.875
$ws.Cells.Item(132, 8).Value = 11369.333  # H132: was 18509.38
$ws.Cells.Item(132, 9).Value = 11954.042  # I132: was 24827.092
$ws.Cells.Item(132, 10).Value = 10199.917  # J132: was 11559.9
$ws.Cells.Item(132, 11).Value = 35862.126  # K132: was 74481.276
$ws.Cells.Item(132, 12).Value = 30599.751  # L132: was 34679.7
$ws.Cells.Item(132, 13).Value = -33332.126  # M132: was -71951.276
$ws.Cells.Item(132, 14).Value = -35659.751  # N132: was -39739.7
$ws.Cells.Item(136, 8).Value = 2208.4792  # H136: was 2434.068
$ws.Cells.Item(136, 9).Value = 1168.6052  # I136: was 1363.6061
$ws.Cells.Item(136, 10).Value = 6160  # J136: was 5645.4546
$ws.Cells.Item(136, 11).Value = 3505.8156  # K136: was 4090.8183
$ws.Cells.Item(136, 12).Value = 18480  # L136: was 16936.3638
$ws.Cells.Item(136, 13).Value = -955.8155999999999  # M136: was -1540.8183
$ws.Cells.Item(136, 14).Value = -23580  # N136: was -22036.3638

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3127.6428  # H122: was 3267.7307
$ws.Cells.Item(122, 9).Value = 2085.7368  # I122: was 2165.3333
$ws.Cells.Item(122, 10).Value = 5327.222  # J122: was 5748.125
$ws.Cells.Item(122, 11).Value = 6257.2104  # K122: was 6495.999899999999
$ws.Cells.Item(122, 12).Value = 15981.666  # L122: was 17244.375
$ws.Cells.Item(122, 13).Value = -3807.2104  # M122: was -4045.999899999999
$ws.Cells.Item(122, 14).Value = -20881.666  # N122: was -22144.375
$ws.Cells.Item(132, 8).Value = 10754473  # H132: was 9525395
$ws.Cells.Item(132, 9).Value = 689.8182  # I132: was 642.37036
$ws.Cells.Item(132, 10).Value = 37041496  # J132: was 41671436
$ws.Cells.Item(132, 11).Value = 2069.4546  # K132: was 1927.11108
$ws.Cells.Item(132, 12).Value = 111124488  # L132: was 125014308
$ws.Cells.Item(132, 13).Value = 460.5454  # M132: was 602.8889199999999
$ws.Cells.Item(132, 14).Value = -111129548  # N132: was -125019368
$ws.Cells.Item(136, 8).Value = 1795.1587  # H136: was 1661.5883
$ws.Cells.Item(136, 9).Value = 506.36957  # I136: was 460.1154
$ws.Cells.Item(136, 10).Value = 5282.4707  # J136: was 5566.375
$ws.Cells.Item(136, 11).Value = 1519.10871  # K136: was 1380.3462
$ws.Cells.Item(136, 12).Value = 15847.4121  # L136: was 16699.125
$ws.Cells.Item(136, 13).Value = 1030.89129  # M136: was 1169.6538
$ws.Cells.Item(136, 14).Value = -20947.4121  # N136: was -21799.125

Write-Host "Updated $($wb.Worksheets.Count) worksheets with refreshed Chocobo profit figures."
